$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 11.08
$ws.Range("D2").Value = 1.69
$ws.Range("F2").Value = 3.85
$ws.Range("K2").Value = 0.77
$ws.Range("N2").Value = 0.92
$ws.Range("Q2").Value = 0.15
$ws.Range("U2").Value = 0.15
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 7
$ws.Range("C4").Value = 45.8
$ws.Range("D4").Value = 13.69
$ws.Range("F4").Value = 1.82
$ws.Range("M4").Value = 0.18
$ws.Range("N4").Value = 0.36
$ws.Range("P4").Value = 0.18
$ws.Range("C5").Value = 21.5
$ws.Range("D5").Value = 16.75
$ws.Range("F5").Value = 6.25
$ws.Range("G5").Value = 0.75
$ws.Range("K5").Value = 0.5
$ws.Range("M5").Value = 0.5
$ws.Range("N5").Value = 0.25
$ws.Range("S5").Value = 0.25
$ws.Range("C6").Value = 62.76
$ws.Range("G6").Value = 0.26
$ws.Range("K6").Value = 0.13
$ws.Range("P6").Value = 0.26
$ws.Range("U6").Value = 0.13
$ws.Range("C7").Value = 20.09
$ws.Range("P7").Value = 0.32
$ws.Range("C8").Value = 61.48
$ws.Range("D8").Value = 26.67
$ws.Range("F8").Value = 7.41
$ws.Range("K8").Value = 0.37
$ws.Range("N8").Value = 0.37
$ws.Range("P8").Value = 0.74
$ws.Range("C9").Value = 37.74
$ws.Range("D9").Value = 44.57
$ws.Range("F9").Value = 10.61
$ws.Range("I9").Value = 0.07000000000000001
$ws.Range("N9").Value = 0.35
$ws.Range("R9").Value = 0.04
$ws.Range("U9").Value = 0.07000000000000001
$ws.Range("C10").Value = 49.52
$ws.Range("D10").Value = 9.470000000000001
$ws.Range("F10").Value = 28.41
$ws.Range("G10").Value = 0.16
$ws.Range("I10").Value = 0.08
$ws.Range("K10").Value = 0.48
$ws.Range("L10").Value = 0.08
$ws.Range("P10").Value = 0.32
$ws.Range("C11").Value = 39.26
$ws.Range("D11").Value = 10.74
$ws.Range("K11").Value = 0.37
$ws.Range("Q11").Value = 1.48
$ws.Range("V11").Value = 1.11
